# Ridership run on 20161026: append data for 24-30 Sep 2016 and refresh
# the existing Sep data (C = Riders, D = Average) through the end of Sep.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# --- Refresh Riders (C) and Average (D) for the existing rows 2-23 ---
$ws.Range("C2").Value  = 193
$ws.Range("D2").Value  = 105.94
$ws.Range("C3").Value  = 149
$ws.Range("D3").Value  = 102.76
$ws.Range("C4").Value  = 97
$ws.Range("D4").Value  = 48.88
$ws.Range("C5").Value  = 77
$ws.Range("D5").Value  = 37.45
$ws.Range("C6").Value  = 183
$ws.Range("D6").Value  = 104.09
$ws.Range("C7").Value  = 240
$ws.Range("D7").Value  = 109.64
$ws.Range("C8").Value  = 227
$ws.Range("D8").Value  = 108.23
$ws.Range("C9").Value  = 210
$ws.Range("D9").Value  = 104.9
$ws.Range("C10").Value = 85
$ws.Range("D10").Value = 49.59
$ws.Range("C11").Value = 73
$ws.Range("D11").Value = 38.13
$ws.Range("C12").Value = 172
$ws.Range("D12").Value = 101.75
$ws.Range("C13").Value = 258
$ws.Range("D13").Value = 106.89
$ws.Range("C14").Value = 240
$ws.Range("D14").Value = 112.06
$ws.Range("C15").Value = 220
$ws.Range("D15").Value = 110.3
$ws.Range("C16").Value = 212
$ws.Range("D16").Value = 107
$ws.Range("C17").Value = 60
$ws.Range("D17").Value = 49.79
$ws.Range("C18").Value = 74
$ws.Range("D18").Value = 38.81
$ws.Range("C19").Value = 229
$ws.Range("D19").Value = 104.35
$ws.Range("C20").Value = 219
$ws.Range("D20").Value = 108.89
$ws.Range("C21").Value = 222
$ws.Range("D21").Value = 114.05
$ws.Range("C22").Value = 191
$ws.Range("D22").Value = 111.76
$ws.Range("C23").Value = 260
$ws.Range("D23").Value = 109.94

# --- Append the new rows of daily data for 24-30 Sep 2016 ---
$ws.Range("A24").Value = "Saturday"
$ws.Range("B24").Value = "24 Sep 2016"
$ws.Range("C24").Value = 107
$ws.Range("D24").Value = 50.87
$ws.Range("E24").Value = 82.36

$ws.Range("A25").Value = "Sunday"
$ws.Range("B25").Value = "25 Sep 2016"
$ws.Range("C25").Value = 67
$ws.Range("D25").Value = 39.33
$ws.Range("E25").Value = 82.53

$ws.Range("A26").Value = "Monday"
$ws.Range("B26").Value = "26 Sep 2016"
$ws.Range("C26").Value = 223
$ws.Range("D26").Value = 106.72
$ws.Range("E26").Value = 82.69

$ws.Range("A27").Value = "Tuesday"
$ws.Range("B27").Value = "27 Sep 2016"
$ws.Range("C27").Value = 189
$ws.Range("D27").Value = 110.3
$ws.Range("E27").Value = 82.85

$ws.Range("A28").Value = "Wednesday"
$ws.Range("B28").Value = "28 Sep 2016"
$ws.Range("C28").Value = 178
$ws.Range("D28").Value = 115.2
$ws.Range("E28").Value = 83.01

$ws.Range("A29").Value = "Thursday"
$ws.Range("B29").Value = "29 Sep 2016"
$ws.Range("C29").Value = 199
$ws.Range("D29").Value = 113.32
$ws.Range("E29").Value = 83.17

$ws.Range("A30").Value = "Friday"
$ws.Range("B30").Value = "30 Sep 2016"
$ws.Range("C30").Value = 266
$ws.Range("D30").Value = 112.89
$ws.Range("E30").Value = 83.33

# --- Extend the three chart series (Ridership, Average, Pilot Target) so
#     they cover the new rows through row 30 instead of row 23 ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = '=SERIES("Ridership",Ridership!$B$2:$B$30,Ridership!$C$2:$C$30,1)'

$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = '=SERIES("Average",Ridership!$B$2:$B$30,Ridership!$D$2:$D$30,2)'

$s3 = $chart.SeriesCollection().Item(3)
$s3.Formula = '=SERIES("Pilot Target",Ridership!$B$2:$B$30,Ridership!$E$2:$E$30,3)'

# --- Move the chart down on the sheet to make room for the extra rows of
#     data (was anchored starting at row 26, now starts at row 33) ---
$co.Top = 487.5
